$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 162.399297
$ws.Cells.Item(2, 8).Value = 487.197891
$ws.Cells.Item(2, 9).Value = 0.3910371682630009
$ws.Cells.Item(2, 10).Value = 0.3910371682630009
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 2.701354
$ws.Cells.Item(2, 14).Value = 8.104061999999999
$ws.Cells.Item(2, 15).Value = 0.02221077311549548
$ws.Cells.Item(2, 16).Value = 0.02221077311549548
$ws.Cells.Item(2, 17).Value = 438.697990548138
$ws.Cells.Item(2, 18).Value = 3948.281914933241
$ws.Cells.Item(2, 19).Value = 0.008685237824015344
$ws.Cells.Item(2, 20).Value = 0.008685237824015344

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 162.399297
$ws.Cells.Item(3, 8).Value = 487.197891
$ws.Cells.Item(3, 9).Value = 0.3910371682630009
$ws.Cells.Item(3, 10).Value = 0.3910371682630009
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 88.14978533333333
$ws.Cells.Item(3, 14).Value = 264.449356
$ws.Cells.Item(3, 15).Value = 0.7247753838328104
$ws.Cells.Item(3, 16).Value = 0.7247753838328105
$ws.Cells.Item(3, 17).Value = 14315.46316883424
$ws.Cells.Item(3, 18).Value = 128839.1685195082
$ws.Cells.Item(3, 19).Value = 0.2834141137207118
$ws.Cells.Item(3, 20).Value = 0.2834141137207118

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 162.399297
$ws.Cells.Item(4, 8).Value = 487.197891
$ws.Cells.Item(4, 9).Value = 0.3910371682630009
$ws.Cells.Item(4, 10).Value = 0.3910371682630009
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 0.24063
$ws.Cells.Item(4, 14).Value = 0.72189
$ws.Cells.Item(4, 15).Value = 0.001978481285600361
$ws.Cells.Item(4, 16).Value = 0.001978481285600361
$ws.Cells.Item(4, 17).Value = 39.07814283711
$ws.Cells.Item(4, 18).Value = 351.70328553399
$ws.Cells.Item(4, 19).Value = 0.0007736597193825069
$ws.Cells.Item(4, 20).Value = 0.0007736597193825069

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 162.399297
$ws.Cells.Item(5, 8).Value = 487.197891
$ws.Cells.Item(5, 9).Value = 0.3910371682630009
$ws.Cells.Item(5, 10).Value = 0.3910371682630009
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 30.53182233333333
$ws.Cells.Item(5, 14).Value = 91.595467
$ws.Cells.Item(5, 15).Value = 0.2510353617660938
$ws.Cells.Item(5, 16).Value = 0.2510353617660938
$ws.Cells.Item(5, 17).Value = 4958.346483062232
$ws.Cells.Item(5, 18).Value = 44625.11834756009
$ws.Cells.Item(5, 19).Value = 0.09816415699889133
$ws.Cells.Item(5, 20).Value = 0.09816415699889133

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 65.41736466666667
$ws.Cells.Item(6, 8).Value = 196.252094
$ws.Cells.Item(6, 9).Value = 0.1575168212364948
$ws.Cells.Item(6, 10).Value = 0.1575168212364948
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 2.701354
$ws.Cells.Item(6, 14).Value = 8.104061999999999
$ws.Cells.Item(6, 15).Value = 0.02221077311549548
$ws.Cells.Item(6, 16).Value = 0.02221077311549548
$ws.Cells.Item(6, 17).Value = 176.7154597117587
$ws.Cells.Item(6, 18).Value = 1590.439137405828
$ws.Cells.Item(6, 19).Value = 0.003498570378357847
$ws.Cells.Item(6, 20).Value = 0.003498570378357847

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 65.41736466666667
$ws.Cells.Item(7, 8).Value = 196.252094
$ws.Cells.Item(7, 9).Value = 0.1575168212364948
$ws.Cells.Item(7, 10).Value = 0.1575168212364948
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 88.14978533333333
$ws.Cells.Item(7, 14).Value = 264.449356
$ws.Cells.Item(7, 15).Value = 0.7247753838328104
$ws.Cells.Item(7, 16).Value = 0.7247753838328105
$ws.Cells.Item(7, 17).Value = 5766.526652439052
$ws.Cells.Item(7, 18).Value = 51898.73987195146
$ws.Cells.Item(7, 19).Value = 0.1141643145718047
$ws.Cells.Item(7, 20).Value = 0.1141643145718047

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 65.41736466666667
$ws.Cells.Item(8, 8).Value = 196.252094
$ws.Cells.Item(8, 9).Value = 0.1575168212364948
$ws.Cells.Item(8, 10).Value = 0.1575168212364948
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 0.24063
$ws.Cells.Item(8, 14).Value = 0.72189
$ws.Cells.Item(8, 15).Value = 0.001978481285600361
$ws.Cells.Item(8, 16).Value = 0.001978481285600361
$ws.Cells.Item(8, 17).Value = 15.74138045974
$ws.Cells.Item(8, 18).Value = 141.67242413766
$ws.Cells.Item(8, 19).Value = 0.0003116440829836626
$ws.Cells.Item(8, 20).Value = 0.0003116440829836626

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 65.41736466666667
$ws.Cells.Item(9, 8).Value = 196.252094
$ws.Cells.Item(9, 9).Value = 0.1575168212364948
$ws.Cells.Item(9, 10).Value = 0.1575168212364948
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 30.53182233333333
$ws.Cells.Item(9, 14).Value = 91.595467
$ws.Cells.Item(9, 15).Value = 0.2510353617660938
$ws.Cells.Item(9, 16).Value = 0.2510353617660938
$ws.Cells.Item(9, 17).Value = 1997.311355517544
$ws.Cells.Item(9, 18).Value = 17975.8021996579
$ws.Cells.Item(9, 19).Value = 0.03954229220334861
$ws.Cells.Item(9, 20).Value = 0.03954229220334861

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 126.3069433333333
$ws.Cells.Item(10, 8).Value = 378.92083
$ws.Cells.Item(10, 9).Value = 0.3041313008456065
$ws.Cells.Item(10, 10).Value = 0.3041313008456065
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 2.701354
$ws.Cells.Item(10, 14).Value = 8.104061999999999
$ws.Cells.Item(10, 15).Value = 0.02221077311549548
$ws.Cells.Item(10, 16).Value = 0.02221077311549548
$ws.Cells.Item(10, 17).Value = 341.1997666012733
$ws.Cells.Item(10, 18).Value = 3070.79789941146
$ws.Cells.Item(10, 19).Value = 0.006754991320402264
$ws.Cells.Item(10, 20).Value = 0.006754991320402266

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 126.3069433333333
$ws.Cells.Item(11, 8).Value = 378.92083
$ws.Cells.Item(11, 9).Value = 0.3041313008456065
$ws.Cells.Item(11, 10).Value = 0.3041313008456065
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 88.14978533333333
$ws.Cells.Item(11, 14).Value = 264.449356
$ws.Cells.Item(11, 15).Value = 0.7247753838328104
$ws.Cells.Item(11, 16).Value = 0.7247753838328105
$ws.Cells.Item(11, 17).Value = 11133.92994094283
$ws.Cells.Item(11, 18).Value = 100205.3694684855
$ws.Cells.Item(11, 19).Value = 0.2204268803059463
$ws.Cells.Item(11, 20).Value = 0.2204268803059464

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 126.3069433333333
$ws.Cells.Item(12, 8).Value = 378.92083
$ws.Cells.Item(12, 9).Value = 0.3041313008456065
$ws.Cells.Item(12, 10).Value = 0.3041313008456065
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 0.24063
$ws.Cells.Item(12, 14).Value = 0.72189
$ws.Cells.Item(12, 15).Value = 0.001978481285600361
$ws.Cells.Item(12, 16).Value = 0.001978481285600361
$ws.Cells.Item(12, 17).Value = 30.3932397743
$ws.Cells.Item(12, 18).Value = 273.5391579687
$ws.Cells.Item(12, 19).Value = 0.0006017180870883257
$ws.Cells.Item(12, 20).Value = 0.0006017180870883258

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 126.3069433333333
$ws.Cells.Item(13, 8).Value = 378.92083
$ws.Cells.Item(13, 9).Value = 0.3041313008456065
$ws.Cells.Item(13, 10).Value = 0.3041313008456065
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 30.53182233333333
$ws.Cells.Item(13, 14).Value = 91.595467
$ws.Cells.Item(13, 15).Value = 0.2510353617660938
$ws.Cells.Item(13, 16).Value = 0.2510353617660938
$ws.Cells.Item(13, 17).Value = 3856.381153319734
$ws.Cells.Item(13, 18).Value = 34707.43037987761
$ws.Cells.Item(13, 19).Value = 0.07634771113216954
$ws.Cells.Item(13, 20).Value = 0.07634771113216955

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 61.180387
$ws.Cells.Item(14, 8).Value = 183.541161
$ws.Cells.Item(14, 9).Value = 0.1473147096548978
$ws.Cells.Item(14, 10).Value = 0.1473147096548978
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 2.701354
$ws.Cells.Item(14, 14).Value = 8.104061999999999
$ws.Cells.Item(14, 15).Value = 0.02221077311549548
$ws.Cells.Item(14, 16).Value = 0.02221077311549548
$ws.Cells.Item(14, 17).Value = 165.269883143998
$ws.Cells.Item(14, 18).Value = 1487.428948295982
$ws.Cells.Item(14, 19).Value = 0.003271973592720026
$ws.Cells.Item(14, 20).Value = 0.003271973592720027

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 61.180387
$ws.Cells.Item(15, 8).Value = 183.541161
$ws.Cells.Item(15, 9).Value = 0.1473147096548978
$ws.Cells.Item(15, 10).Value = 0.1473147096548978
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 88.14978533333333
$ws.Cells.Item(15, 14).Value = 264.449356
$ws.Cells.Item(15, 15).Value = 0.7247753838328104
$ws.Cells.Item(15, 16).Value = 0.7247753838328105
$ws.Cells.Item(15, 17).Value = 5393.037980660257
$ws.Cells.Item(15, 18).Value = 48537.34182594231
$ws.Cells.Item(15, 19).Value = 0.1067700752343476
$ws.Cells.Item(15, 20).Value = 0.1067700752343476

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 61.180387
$ws.Cells.Item(16, 8).Value = 183.541161
$ws.Cells.Item(16, 9).Value = 0.1473147096548978
$ws.Cells.Item(16, 10).Value = 0.1473147096548978
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 0.24063
$ws.Cells.Item(16, 14).Value = 0.72189
$ws.Cells.Item(16, 15).Value = 0.001978481285600361
$ws.Cells.Item(16, 16).Value = 0.001978481285600361
$ws.Cells.Item(16, 17).Value = 14.72183652381
$ws.Cells.Item(16, 18).Value = 132.49652871429
$ws.Cells.Item(16, 19).Value = 0.0002914593961458662
$ws.Cells.Item(16, 20).Value = 0.0002914593961458662

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 61.180387
$ws.Cells.Item(17, 8).Value = 183.541161
$ws.Cells.Item(17, 9).Value = 0.1473147096548978
$ws.Cells.Item(17, 10).Value = 0.1473147096548978
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 30.53182233333333
$ws.Cells.Item(17, 14).Value = 91.595467
$ws.Cells.Item(17, 15).Value = 0.2510353617660938
$ws.Cells.Item(17, 16).Value = 0.2510353617660938
$ws.Cells.Item(17, 17).Value = 1867.948706168576
$ws.Cells.Item(17, 18).Value = 16811.53835551719
$ws.Cells.Item(17, 19).Value = 0.03698120143168434
$ws.Cells.Item(17, 20).Value = 0.03698120143168435
